$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.268.20'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '1.657.33'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  -0.75%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.28'
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5245'
$ws.Range('E6').Value = '  -1.67%  '
$ws.Range('E7').Value = '  -0.65%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2672'
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06372'
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.72'
$ws.Range('E10').Value = '  -0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07716'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.605'
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('D13').Value = '1.619.14'
$ws.Range('E13').Value = '  -3.34%  '
$ws.Range('D14').Value = '1.885.54'
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5653'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').Value = '0.0₅8263'
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.50'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').Value = '26.261.85'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.696'
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.29'
$ws.Range('E21').Value = '  -2.42%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.44'
$ws.Range('E22').Value = '  +1.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.009'
$ws.Range('E23').Value = '  -1.04%  '
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.32'
$ws.Range('E25').Value = '  -2.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1203'
$ws.Range('E26').Value = '  -1.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.299'
$ws.Range('E27').Value = '  +0.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.94'
$ws.Range('E28').Value = '  -1.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.509'
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05657'
$ws.Range('E30').Value = '  -4.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.277'
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.512'
$ws.Range('E32').Value = '  -1.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.355'
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.585'
$ws.Range('E34').Value = '  -1.58%  '
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9493'
$ws.Range('E36').Value = '  -2.00%  '
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5777'
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01603'
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.977'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.568'
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8464'
$ws.Range('E42').Value = '  -2.14%  '
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.90'
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('D45').Value = '1.018.45'
$ws.Range('E45').Value = '  -5.55%  '
$ws.Range('D46').Value = '1.796.41'
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '58.39'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('D48').Value = '0.0₈106'
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('E49').Value = '  -1.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05327'
$ws.Range('E50').Value = '  +3.16%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4349'
$ws.Range('E51').Value = '  -1.63%  '
